# Rename the sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Sheet {{stringName}} ABC"

# Re-write the rich-text cell A3 with the new runs/formatting
$cellA3 = $ws.Range("A3")
$cellA3.Value = "Even with Text formatting: Some bold {{stringName}}! But only the cell font will be preserved. Any text-section specific formatting will be discarded."

$cellA3.Characters(99,13).Font.Bold = $true
$cellA3.Characters(99,13).Font.Italic = $true

$cellA3.Characters(113,8).Font.Bold = $true
$cellA3.Characters(113,8).Font.Underline = $true

$cellA3.Characters(121,1).Font.Bold = $true

$cellA3.Characters(122,10).Font.Italic = $true
$cellA3.Characters(122,10).Font.Underline = $true

$cellA3.Characters(132,19).Font.Bold = $true

$cellA3.Font.Bold = $true

# New cells
$ws.Range("A7").Value = "Empty lines are not a problem."

$cellF9 = $ws.Range("F9")
$cellF9.Value = "{{stringName}} folks, Let Jebt take on the {{stringName}}"
$cellF9.Font.Italic = $true
$cellF9.Font.Underline = $true

# New selection, matching target sheetView
$ws.Range("G5").Select() | Out-Null
